$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric need to be force-typed as text
# (matching the original inline-string cell type) before assignment,
# otherwise Excel auto-converts them to numbers and drops formatting
# such as trailing zeros (e.g. "1.40" -> 1.4). NumberFormat is applied
# per-cell (not as a multi-area union) and reset to Normal afterwards so
# no visible style change is left on the cell.
$textForceCells = @("D13", "D16", "D20", "D22", "D27", "D28", "D29", "D30", "D33", "D36", "D41", "D44", "D49")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D13").Value = "14.83"
$ws.Range("D16").Value = "5.33"
$ws.Range("D20").Value = "71.52"
$ws.Range("D22").Value = "229.13"
$ws.Range("D27").Value = "171.96"
$ws.Range("D28").Value = "0.137"
$ws.Range("D29").Value = "1.40"
$ws.Range("D30").Value = "19.42"
$ws.Range("D33").Value = "0.0630"
$ws.Range("D36").Value = "1.82"
$ws.Range("D41").Value = "101.65"
$ws.Range("D44").Value = "16.96"
$ws.Range("D49").Value = "7.34"

foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}

# Remaining cell updates (values that do not look numeric, so Excel keeps
# them as text automatically)
$ws.Range("D2").Value = "37.793.71"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").Value = "2.076.61"
$ws.Range("E3").Value = "  -0.56%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  -0.65%  "
$ws.Range("E7").Value = "  -0.89%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("E10").Value = "  -1.10%  "
$ws.Range("E11").Value = "  +2.86%  "
$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").Value = "2.384.03"
$ws.Range("E12").Value = "  -0.50%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("E14").Value = "  -1.65%  "
$ws.Range("E15").Value = "  +1.16%  "
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("D17").Value = "2.065.78"
$ws.Range("E17").Value = "  -0.98%  "
$ws.Range("D18").Value = "37.748.53"
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("E19").Value = "  -1.90%  "
$ws.Range("E20").Value = "  -0.10%  "
$ws.Range("D21").Value = "0.0₃0841"
$ws.Range("E21").Value = "  +1.20%  "
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("E24").Value = "  -1.04%  "
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("E26").Value = "  +6.52%  "
$ws.Range("E27").Value = "  +0.73%  "
$ws.Range("E28").Value = "  -0.94%  "
$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("E29").Value = "  -1.54%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("E30").Value = "  -0.80%  "
$ws.Range("E31").Value = "  +0.89%  "
$ws.Range("E33").Value = "  -0.31%  "
$ws.Range("E34").Value = "  -1.34%  "
$ws.Range("E35").Value = "  -2.30%  "
$ws.Range("E36").Value = "  -0.64%  "
$ws.Range("E37").Value = "  -3.06%  "
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("E39").Value = "  -1.10%  "
$ws.Range("E40").Value = "  +7.00%  "
$ws.Range("E41").Value = "  +2.07%  "
$ws.Range("E42").Value = "  -1.09%  "
$ws.Range("E43").Value = "  -0.19%  "
$ws.Range("E44").Value = "  +4.98%  "
$ws.Range("D45").Value = "1.450.83"
$ws.Range("E45").Value = "  -0.79%  "
$ws.Range("E46").Value = "  -1.90%  "
$ws.Range("E47").Value = "  -1.51%  "
$ws.Range("E48").Value = "  -4.51%  "
$ws.Range("E49").Value = "  -2.08%  "
$ws.Range("E50").Value = "  -1.51%  "
$ws.Range("D51").Value = "2.270.11"
$ws.Range("E51").Value = "  -0.48%  "
